$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Drop the two unused trailing columns (H, I). This also automatically
#    shrinks the dimension and collapses the F2:I2 merge down to F2:G2, and
#    removes the old H8 gantt marker that lived in the deleted column.
# ---------------------------------------------------------------------------
$ws.Range("H1:I1").EntireColumn.Delete()

# ---------------------------------------------------------------------------
# 2. Month / week headers (row 2-3)
# ---------------------------------------------------------------------------
$ws.Range("F2").Value = "January"
$ws.Range("F3").Value = "01/Jan - 07/Jan"
$ws.Range("G3").Value = "08/Jan - 14/Jan"

# ---------------------------------------------------------------------------
# 3. Row 4 becomes a blank spacer row (old "Task 1" / "ML1 - T1" + marker removed)
# ---------------------------------------------------------------------------
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("F4").Clear()

# ---------------------------------------------------------------------------
# 4. Row 5 becomes "Task 1.1" with its own date span + gantt marker in F5
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "Task 1.1"
$ws.Range("C5").Value = "M1 - T1"
$ws.Range("D5").Value = "01/01"
$ws.Range("E5").Value = "01/07"
$ws.Range("F5").Interior.Color = 42495

# ---------------------------------------------------------------------------
# 5. Row 6 becomes "Task 1.2" with its own date span; marker moves F6 -> G6
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = "Task 1.2"
$ws.Range("C6").Value = "M1 - T2"
$ws.Range("D6").Value = "01/08"
$ws.Range("E6").Value = "01/14"
$ws.Range("F6").Clear()
$ws.Range("G6").Interior.Color = 42495

# ---------------------------------------------------------------------------
# 6. Row 7 becomes a blank spacer row (old "Task 3" / "ML2 - X2" + marker removed)
# ---------------------------------------------------------------------------
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("G7").Clear()

# ---------------------------------------------------------------------------
# 7. Row 8 becomes "Task 2.1" with its own date span; marker moves to F8
# ---------------------------------------------------------------------------
$ws.Range("B8").Value = "Task 2.1"
$ws.Range("C8").Value = "M2 - X1"
$ws.Range("D8").Value = "01/15"
$ws.Range("E8").Value = "01/21"
$ws.Range("F8").Interior.Color = 42495

# ---------------------------------------------------------------------------
# 8. New row 9: "Task 2.2" with its own date span + gantt marker in G9
# ---------------------------------------------------------------------------
$ws.Range("B9").WrapText = $true
$ws.Range("B9").Value = "Task 2.2"
$ws.Range("C9").WrapText = $true
$ws.Range("C9").Value = "M2 - X2"
$ws.Range("D9").Value = "01/22"
$ws.Range("E9").Value = "01/28"
$ws.Range("G9").Interior.Color = 42495
